$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("B4").Value = "[André Guimarães-CAD, -]"
$ws.Range("C4").Value = "[-, João Paulo-Sistemas Digitais, -]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[Pedro Bispo-Acionamentos Elétricos, -, Sandro-Programação de Computadores, -]"

$ws.Range("B6").Value = "[André Guimarães-CAD, -]"
$ws.Range("C6").Value = "[-, João Paulo-Sistemas Digitais, -]"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "[Pedro Bispo-Acionamentos Elétricos, -, Sandro-Programação de Computadores, -]"
$ws.Range("F6").Value = "Sandro-Circuitos Elétrico"

$ws.Range("B7").Value = "[André Guimarães-CAD, -]"
$ws.Range("C7").Value = "[-, João Paulo-Sistemas Digitais, -]"
$ws.Range("D7").Value = "[-, -, João Paulo-Sistemas Digitais]"
$ws.Range("E7").Value = "[Pedro Bispo-Acionamentos Elétricos, -, Sandro-Programação de Computadores, -]"

$ws.Range("B8").Value = "[André Guimarães-CAD, -]"
$ws.Range("C8").Value = "[-, João Paulo-Sistemas Digitais, -]"
$ws.Range("D8").Value = "[-, -, João Paulo-Sistemas Digitais]"
$ws.Range("E8").Value = "[Pedro Bispo-Acionamentos Elétricos, -, Sandro-Programação de Computadores, -]"
$ws.Range("F8").Value = "Nilton-M. T. R."

$ws.Range("C11").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("C14").Value = "Cleidson-Eletronica Analógica e de Potên"
$ws.Range("C15").Value = "Cleidson-Eletronica Analógica e de Potên"
